$wb = $excel.ActiveWorkbook

# --- Step1_Data (71 cell updates) ---
$ws = $wb.Worksheets.Item("Step1_Data")
$ws.Range("D2").Value = 0.0008464570636233179
$ws.Range("E2").Value = 0.07009297764934155
$ws.Range("F2").Value = 0.04448455179163968
$ws.Range("G2").Value = 0.04539466964219265
$ws.Range("H2").Value = 0.02037722750575902
$ws.Range("K2").Value = 0.04850738180508758
$ws.Range("M2").Value = 0.1101160843398459
$ws.Range("N2").Value = 0.1412450804082384
$ws.Range("O2").Value = 0.06786809913577249
$ws.Range("P2").Value = 0.03619626360227308
$ws.Range("Q2").Value = 0.01548960142468087
$ws.Range("R2").Value = 0.1332746214499288
$ws.Range("S2").Value = 0.07247048239457612
$ws.Range("T2").Value = 0.1851876904402279
$ws.Range("AC2").Value = 0.005782613102825118
$ws.Range("AG2").Value = 0.002381642793435917
$ws.Range("AH2").Value = 0.0002845554505514655
$ws.Range("D3").Value = 0.171702537436719
$ws.Range("E3").Value = 0.001949443460812957
$ws.Range("F3").Value = 0.05826936416561582
$ws.Range("G3").Value = 0.03383913379315021
$ws.Range("K3").Value = 0.07604134795866979
$ws.Range("L3").Value = 0.08210080009524807
$ws.Range("N3").Value = 0.243783329372892
$ws.Range("Q3").Value = 0.0806443918020239
$ws.Range("R3").Value = 0.01067298868354193
$ws.Range("S3").Value = 0.1706124962805226
$ws.Range("T3").Value = 0.008889269952726517
$ws.Range("U3").Value = 0.03965709155553353
$ws.Range("AB3").Value = 0.02145489129396419
$ws.Range("AC3").Value = 0.0003829141485795281
$ws.Range("D4").Value = 0.2207149227007408
$ws.Range("E4").Value = 0.008132993378711265
$ws.Range("F4").Value = 0.06818760751432139
$ws.Range("G4").Value = 0.05664690838283121
$ws.Range("K4").Value = 0.02750203285179573
$ws.Range("L4").Value = 0.1136939909426669
$ws.Range("N4").Value = 0.1917019104425311
$ws.Range("Q4").Value = 0.08658711917312777
$ws.Range("S4").Value = 0.1922159363556849
$ws.Range("T4").Value = 0.004069317007811368
$ws.Range("U4").Value = 0.01865310843509864
$ws.Range("AB4").Value = 0.01189415281467878
$ws.Range("E5").Value = 0.1710063744167687
$ws.Range("F5").Value = 0.02667952771242
$ws.Range("G5").Value = 0.07048842176554446
$ws.Range("H5").Value = 0.04259321876820902
$ws.Range("L5").Value = 0.04776879471852993
$ws.Range("M5").Value = 0.1148558486469529
$ws.Range("N5").Value = 0.008975720043108143
$ws.Range("O5").Value = 0.1738507635505282
$ws.Range("R5").Value = 0.102983523309786
$ws.Range("T5").Value = 0.1883609295811709
$ws.Range("V5").Value = 0.03305416598677679
$ws.Range("AC5").Value = 0.01938271150020498
$ws.Range("D6").Value = 0.03458282264795465
$ws.Range("E6").Value = 0.02312423813842653
$ws.Range("F6").Value = 0.08316488917773351
$ws.Range("G6").Value = 0.09996811038135137
$ws.Range("H6").Value = 0.02769317061367617
$ws.Range("K6").Value = 0.03538183398599539
$ws.Range("M6").Value = 0.05953123838766115
$ws.Range("N6").Value = 0.1970872778162309
$ws.Range("O6").Value = 0.0618305601271338
$ws.Range("P6").Value = 0.007059615658957061
$ws.Range("Q6").Value = 0.01314662652846956
$ws.Range("R6").Value = 0.1247741826188429
$ws.Range("S6").Value = 0.03666306034347517
$ws.Range("T6").Value = 0.1476565635511564
$ws.Range("U6").Value = 0.01004633560979529
$ws.Range("AC6").Value = 0.03828947441314022

# --- Step2_Sj (145 cell updates) ---
$ws = $wb.Worksheets.Item("Step2_Sj")
$ws.Range("D2").Value = 0.0008464570636233179
$ws.Range("E2").Value = 0.07093943471296488
$ws.Range("F2").Value = 0.1154239865046046
$ws.Range("G2").Value = 0.1608186561467972
$ws.Range("H2").Value = 0.1811958836525562
$ws.Range("I2").Value = 0.1811958836525562
$ws.Range("J2").Value = 0.1811958836525562
$ws.Range("K2").Value = 0.2297032654576438
$ws.Range("L2").Value = 0.2297032654576438
$ws.Range("M2").Value = 0.3398193497974897
$ws.Range("N2").Value = 0.4810644302057281
$ws.Range("O2").Value = 0.5489325293415006
$ws.Range("P2").Value = 0.5851287929437737
$ws.Range("Q2").Value = 0.6006183943684547
$ws.Range("R2").Value = 0.7338930158183834
$ws.Range("S2").Value = 0.8063634982129595
$ws.Range("T2").Value = 0.9915511886531874
$ws.Range("U2").Value = 0.9915511886531874
$ws.Range("V2").Value = 0.9915511886531874
$ws.Range("W2").Value = 0.9915511886531874
$ws.Range("X2").Value = 0.9915511886531874
$ws.Range("Y2").Value = 0.9915511886531874
$ws.Range("Z2").Value = 0.9915511886531874
$ws.Range("AA2").Value = 0.9915511886531874
$ws.Range("AB2").Value = 0.9915511886531874
$ws.Range("AC2").Value = 0.9973338017560125
$ws.Range("AD2").Value = 0.9973338017560125
$ws.Range("AE2").Value = 0.9973338017560125
$ws.Range("AF2").Value = 0.9973338017560125
$ws.Range("AG2").Value = 0.9997154445494484
$ws.Range("D3").Value = 0.171702537436719
$ws.Range("E3").Value = 0.1736519808975319
$ws.Range("F3").Value = 0.2319213450631478
$ws.Range("G3").Value = 0.265760478856298
$ws.Range("H3").Value = 0.265760478856298
$ws.Range("I3").Value = 0.265760478856298
$ws.Range("J3").Value = 0.265760478856298
$ws.Range("K3").Value = 0.3418018268149678
$ws.Range("L3").Value = 0.4239026269102159
$ws.Range("M3").Value = 0.4239026269102159
$ws.Range("N3").Value = 0.6676859562831079
$ws.Range("O3").Value = 0.6676859562831079
$ws.Range("P3").Value = 0.6676859562831079
$ws.Range("Q3").Value = 0.7483303480851318
$ws.Range("R3").Value = 0.7590033367686737
$ws.Range("S3").Value = 0.9296158330491963
$ws.Range("T3").Value = 0.9385051030019228
$ws.Range("U3").Value = 0.9781621945574563
$ws.Range("V3").Value = 0.9781621945574563
$ws.Range("W3").Value = 0.9781621945574563
$ws.Range("X3").Value = 0.9781621945574563
$ws.Range("Y3").Value = 0.9781621945574563
$ws.Range("Z3").Value = 0.9781621945574563
$ws.Range("AA3").Value = 0.9781621945574563
$ws.Range("AB3").Value = 0.9996170858514205
$ws.Range("D4").Value = 0.2207149227007408
$ws.Range("E4").Value = 0.228847916079452
$ws.Range("F4").Value = 0.2970355235937734
$ws.Range("G4").Value = 0.3536824319766046
$ws.Range("H4").Value = 0.3536824319766046
$ws.Range("I4").Value = 0.3536824319766046
$ws.Range("J4").Value = 0.3536824319766046
$ws.Range("K4").Value = 0.3811844648284003
$ws.Range("L4").Value = 0.4948784557710672
$ws.Range("M4").Value = 0.4948784557710672
$ws.Range("N4").Value = 0.6865803662135983
$ws.Range("O4").Value = 0.6865803662135983
$ws.Range("P4").Value = 0.6865803662135983
$ws.Range("Q4").Value = 0.7731674853867261
$ws.Range("R4").Value = 0.7731674853867261
$ws.Range("S4").Value = 0.965383421742411
$ws.Range("T4").Value = 0.9694527387502223
$ws.Range("U4").Value = 0.9881058471853209
$ws.Range("V4").Value = 0.9881058471853209
$ws.Range("W4").Value = 0.9881058471853209
$ws.Range("X4").Value = 0.9881058471853209
$ws.Range("Y4").Value = 0.9881058471853209
$ws.Range("Z4").Value = 0.9881058471853209
$ws.Range("AA4").Value = 0.9881058471853209
$ws.Range("AB4").Value = 0.9999999999999997
$ws.Range("AC4").Value = 0.9999999999999997
$ws.Range("AD4").Value = 0.9999999999999997
$ws.Range("AE4").Value = 0.9999999999999997
$ws.Range("AF4").Value = 0.9999999999999997
$ws.Range("AG4").Value = 0.9999999999999997
$ws.Range("AH4").Value = 0.9999999999999997
$ws.Range("AI4").Value = 0.9999999999999997
$ws.Range("AJ4").Value = 0.9999999999999997
$ws.Range("E5").Value = 0.1710063744167687
$ws.Range("F5").Value = 0.1976859021291887
$ws.Range("G5").Value = 0.2681743238947331
$ws.Range("H5").Value = 0.3107675426629422
$ws.Range("I5").Value = 0.3107675426629422
$ws.Range("J5").Value = 0.3107675426629422
$ws.Range("K5").Value = 0.3107675426629422
$ws.Range("L5").Value = 0.3585363373814721
$ws.Range("M5").Value = 0.473392186028425
$ws.Range("N5").Value = 0.4823679060715331
$ws.Range("O5").Value = 0.6562186696220613
$ws.Range("P5").Value = 0.6562186696220613
$ws.Range("Q5").Value = 0.6562186696220613
$ws.Range("R5").Value = 0.7592021929318473
$ws.Range("S5").Value = 0.7592021929318473
$ws.Range("T5").Value = 0.9475631225130182
$ws.Range("U5").Value = 0.9475631225130182
$ws.Range("V5").Value = 0.980617288499795
$ws.Range("W5").Value = 0.980617288499795
$ws.Range("X5").Value = 0.980617288499795
$ws.Range("Y5").Value = 0.980617288499795
$ws.Range("Z5").Value = 0.980617288499795
$ws.Range("AA5").Value = 0.980617288499795
$ws.Range("AB5").Value = 0.980617288499795
$ws.Range("D6").Value = 0.03458282264795465
$ws.Range("E6").Value = 0.05770706078638118
$ws.Range("F6").Value = 0.1408719499641147
$ws.Range("G6").Value = 0.2408400603454661
$ws.Range("H6").Value = 0.2685332309591422
$ws.Range("I6").Value = 0.2685332309591422
$ws.Range("J6").Value = 0.2685332309591422
$ws.Range("K6").Value = 0.3039150649451376
$ws.Range("L6").Value = 0.3039150649451376
$ws.Range("M6").Value = 0.3634463033327988
$ws.Range("N6").Value = 0.5605335811490296
$ws.Range("O6").Value = 0.6223641412761635
$ws.Range("P6").Value = 0.6294237569351205
$ws.Range("Q6").Value = 0.6425703834635901
$ws.Range("R6").Value = 0.7673445660824331
$ws.Range("S6").Value = 0.8040076264259083
$ws.Range("T6").Value = 0.9516641899770647
$ws.Range("U6").Value = 0.96171052558686
$ws.Range("V6").Value = 0.96171052558686
$ws.Range("W6").Value = 0.96171052558686
$ws.Range("X6").Value = 0.96171052558686
$ws.Range("Y6").Value = 0.96171052558686
$ws.Range("Z6").Value = 0.96171052558686
$ws.Range("AA6").Value = 0.96171052558686
$ws.Range("AB6").Value = 0.96171052558686
$ws.Range("AC6").Value = 1
$ws.Range("AD6").Value = 1
$ws.Range("AE6").Value = 1
$ws.Range("AF6").Value = 1
$ws.Range("AG6").Value = 1
$ws.Range("AH6").Value = 1
$ws.Range("AI6").Value = 1
$ws.Range("AJ6").Value = 1

# --- Step3_DataPts_0.5 (8 cell updates) ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.5")
$ws.Range("F2").Value = 0.5489325293415006
$ws.Range("F3").Value = 0.6676859562831079
$ws.Range("F4").Value = 0.6865803662135983
$ws.Range("F5").Value = 0.6562186696220613
$ws.Range("C6").Value = 2
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0.5605335811490296
$ws.Range("G6").Value = 11

# --- Step3_DataPts_0.7 (12 cell updates) ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.7")
$ws.Range("D2").Value = 17
$ws.Range("F2").Value = 0.7338930158183834
$ws.Range("G2").Value = 15
$ws.Range("D3").Value = 16
$ws.Range("F3").Value = 0.7483303480851318
$ws.Range("G3").Value = 15
$ws.Range("F4").Value = 0.7731674853867261
$ws.Range("F5").Value = 0.7592021929318473
$ws.Range("C6").Value = 2
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0.7673445660824331
$ws.Range("G6").Value = 15

# --- Step3_DataPts_0.8 (10 cell updates) ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.8")
$ws.Range("D2").Value = 18
$ws.Range("F2").Value = 0.8063634982129595
$ws.Range("G2").Value = 16
$ws.Range("F3").Value = 0.9296158330491963
$ws.Range("F4").Value = 0.965383421742411
$ws.Range("F5").Value = 0.9475631225130182
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 18
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0.8040076264259083

# --- Step3_DataPts_0.9 (8 cell updates) ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.9")
$ws.Range("F2").Value = 0.9915511886531874
$ws.Range("F3").Value = 0.9296158330491963
$ws.Range("F4").Value = 0.965383421742411
$ws.Range("F5").Value = 0.9475631225130182
$ws.Range("C6").Value = 2
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0.9516641899770647
$ws.Range("G6").Value = 17
